# feat: add 2022-Q1 data
#
# Before: sheets are "2021-Q1", "总计".
# After:  sheets are "2021-Q1", "2022-Q1", "总计" (总计 updated with a new
#         2022-Q1 row).
#
# Strategy:
#   1. Repurpose the existing "总计" worksheet (2nd sheet) into "2022-Q1":
#      it already carries the bold/bordered header style (s=2) on B1:D1
#      and the matching index style on A2, so we reuse those cells and
#      extend the sheet with new header/data columns and rows.
#   2. Append a brand-new "总计" worksheet right after "2022-Q1", and fill
#      it with the summary table (2022-Q1 first, 2021-Q1 second).
#
# Notes on cell typing: several numeric-looking fields (fund code, fund
# scale, position %, etc.) must be stored as literal text, not numbers
# (matching how the sheet already stores "009263", "0.66", etc. on the
# "2021-Q1" sheet). Assigning a numeric-looking string straight to
# .Value auto-converts it to a number (and would strip the leading
# zeros from fund codes like "009263"), so for those cells we briefly
# mark the range as Text (NumberFormat "@") before assigning the value,
# then clear the formatting again afterwards so the cell ends up with
# no explicit style, exactly like the sheet's other plain data cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: turn the existing "总计" worksheet into "2022-Q1"
# ---------------------------------------------------------------------
$ws2022 = $wb.Worksheets.Item(2)
$ws2022.Name = "2022-Q1"

# --- Header row -------------------------------------------------------
# B1:D1 already have the bold/centered/bordered style; just change text.
$ws2022.Range("B1").Value = "基金代码"
$ws2022.Range("C1").Value = "基金名称"
$ws2022.Range("D1").Value = "基金规模"

# E1:H1 are new cells - copy the existing header's formatting across so
# they pick up the same style rather than the default one.
$ws2022.Range("B1").Copy() | Out-Null
$ws2022.Range("E1:H1").PasteSpecial(-4122) | Out-Null
$ws2022.Range("E1").Value = "股票总仓位"
$ws2022.Range("F1").Value = "仓位占比"
$ws2022.Range("G1").Value = "持有市值(亿元)"
$ws2022.Range("H1").Value = "仓位排名"

# --- Index column (A) --------------------------------------------------
# A2 already carries the index style; copy it down onto the new rows.
$ws2022.Range("A2").Copy() | Out-Null
$ws2022.Range("A3:A5").PasteSpecial(-4122) | Out-Null
$ws2022.Range("A2").Value = 0
$ws2022.Range("A3").Value = 1
$ws2022.Range("A4").Value = 2
$ws2022.Range("A5").Value = 3

# --- Data rows (B:G as text, H as a plain number) ----------------------
$textRange = $ws2022.Range("B2:G5")
$textRange.NumberFormat = "@"

$ws2022.Range("B2").Value = "519625"
$ws2022.Range("C2").Value = "银河君盛灵活配置混合A"
$ws2022.Range("D2").Value = "4.91"
$ws2022.Range("E2").Value = "20.05"
$ws2022.Range("F2").Value = "1.74"
$ws2022.Range("G2").Value = "0.0854"

$ws2022.Range("B3").Value = "519626"
$ws2022.Range("C3").Value = "银河君盛灵活配置混合C"
$ws2022.Range("D3").Value = "2.33"
$ws2022.Range("E3").Value = "20.05"
$ws2022.Range("F3").Value = "1.74"
$ws2022.Range("G3").Value = "0.0405"

$ws2022.Range("B4").Value = "009263"
$ws2022.Range("C4").Value = "华宝红利精选混合A"
$ws2022.Range("D4").Value = "0.46"
$ws2022.Range("E4").Value = "83.67"
$ws2022.Range("F4").Value = "0.96"
$ws2022.Range("G4").Value = "0.0044"

$ws2022.Range("B5").Value = "010841"
$ws2022.Range("C5").Value = "华宝红利精选混合C"
$ws2022.Range("D5").Value = "0.16"
$ws2022.Range("E5").Value = "83.67"
$ws2022.Range("F5").Value = "0.96"
$ws2022.Range("G5").Value = "0.0015"

# Drop the temporary "@" formatting now that the literal text values are
# set, so these cells end up with no explicit style (matching the rest
# of the sheet's plain data cells).
$textRange.ClearFormats()

$ws2022.Range("H2").Value = 4
$ws2022.Range("H3").Value = 4
$ws2022.Range("H4").Value = 10
$ws2022.Range("H5").Value = 10

# ---------------------------------------------------------------------
# Step 2: add a brand-new "总计" worksheet right after "2022-Q1"
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Add($null, $ws2022)
$wsTotal.Name = "总计"

# Borrow the header/index styling from "2022-Q1" so the new sheet's
# formatted cells (bold, centered, thin border) match.
$ws2022.Range("B1").Copy() | Out-Null
$wsTotal.Range("B1:D1").PasteSpecial(-4122) | Out-Null
$ws2022.Range("A2").Copy() | Out-Null
$wsTotal.Range("A2:A3").PasteSpecial(-4122) | Out-Null

$wsTotal.Range("B1").Value = "日期"
$wsTotal.Range("C1").Value = "持有数量(只)"
$wsTotal.Range("D1").Value = "持有市值(亿元)"

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 4
$wsTotal.Range("D2").Value = 0.13

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q1"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.01
